# Auto-generated Excel COM-interop script to apply the scheduled market-data refresh
# to the per-job "Leve Profit" sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 649.5263
$ws.Range("I28").Value = 546.3125
$ws.Range("J28").Value = 1200
$ws.Range("K28").Value = 546.3125
$ws.Range("L28").Value = 1200
$ws.Range("M28").Value = -61.3125
$ws.Range("N28").Value = -2170
$ws.Range("H98").Value = 2999.7896
$ws.Range("I98").Value = 1374.75
$ws.Range("K98").Value = 1374.75
$ws.Range("M98").Value = 123.25
$ws.Range("H122").Value = 2999.7896
$ws.Range("I122").Value = 1374.75
$ws.Range("K122").Value = 4124.25
$ws.Range("M122").Value = -1674.25
$ws.Range("H132").Value = 921216.75
$ws.Range("I132").Value = 14672.25
$ws.Range("K132").Value = 44016.75
$ws.Range("M132").Value = -41486.75
$ws.Range("H137").Value = 2351.9656
$ws.Range("I137").Value = 1621.9
$ws.Range("J137").Value = 3974.3333
$ws.Range("K137").Value = 4865.700000000001
$ws.Range("L137").Value = 11922.9999
$ws.Range("M137").Value = -2315.700000000001
$ws.Range("N137").Value = -17022.9999
$ws.Range("H138").Value = 3650.54
$ws.Range("I138").Value = 521.3
$ws.Range("J138").Value = 4991.643
$ws.Range("K138").Value = 1563.9
$ws.Range("L138").Value = 14974.929
$ws.Range("M138").Value = 3576.1
$ws.Range("N138").Value = -25254.929
$ws.Range("H141").Value = 6479.316
$ws.Range("I141").Value = 7567.6333
$ws.Range("J141").Value = 2398.125
$ws.Range("K141").Value = 22702.8999
$ws.Range("L141").Value = 7194.375
$ws.Range("M141").Value = -17522.8999
$ws.Range("N141").Value = -17554.375

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3642.2056
$ws.Range("I32").Value = 3702.6604
$ws.Range("J32").Value = 3482
$ws.Range("K32").Value = 3702.6604
$ws.Range("L32").Value = 3482
$ws.Range("M32").Value = -3415.6604
$ws.Range("N32").Value = -4056
$ws.Range("H38").Value = 8000
$ws.Range("I38").Value = 8000
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 8000
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -7533
$ws.Range("N38").ClearContents()
$ws.Range("H64").Value = 36000
$ws.Range("J64").Value = 36000
$ws.Range("L64").Value = 36000
$ws.Range("N64").Value = -36496
$ws.Range("H67").Value = 36000
$ws.Range("J67").Value = 36000
$ws.Range("L67").Value = 36000
$ws.Range("N67").Value = -37716

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H134").Value = 2350.9714
$ws.Range("I134").Value = 1536.8
$ws.Range("J134").Value = 4386.4
$ws.Range("K134").Value = 4610.4
$ws.Range("L134").Value = 13159.2
$ws.Range("M134").Value = -2075.4
$ws.Range("N134").Value = -18229.2

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22730462
$ws.Range("I31").Value = 1220.25
$ws.Range("K31").Value = 1220.25
$ws.Range("M31").Value = -925.25
$ws.Range("H34").Value = 22730462
$ws.Range("I34").Value = 1220.25
$ws.Range("K34").Value = 1220.25
$ws.Range("M34").Value = -1018.25
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H39").Value = 20463.5
$ws.Range("I39").Value = 3074
$ws.Range("J39").Value = 25206.092
$ws.Range("K39").Value = 3074
$ws.Range("L39").Value = 25206.092
$ws.Range("M39").Value = -2683
$ws.Range("N39").Value = -25988.092
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H49").Value = 20463.5
$ws.Range("I49").Value = 3074
$ws.Range("J49").Value = 25206.092
$ws.Range("K49").Value = 3074
$ws.Range("L49").Value = 25206.092
$ws.Range("M49").Value = -2892
$ws.Range("N49").Value = -25570.092
$ws.Range("H134").Value = 6462.227
$ws.Range("I134").Value = 9028.166999999999
$ws.Range("J134").Value = 3383.1
$ws.Range("K134").Value = 27084.501
$ws.Range("L134").Value = 10149.3
$ws.Range("M134").Value = -24549.501
$ws.Range("N134").Value = -15219.3

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 396.45456
$ws.Range("J92").Value = 393.5
$ws.Range("L92").Value = 1180.5
$ws.Range("N92").Value = -3676.5
$ws.Range("H113").Value = 790.5417
$ws.Range("I113").Value = 780.6842
$ws.Range("J113").Value = 828
$ws.Range("K113").Value = 2342.0526
$ws.Range("L113").Value = 2484
$ws.Range("M113").Value = -172.0526
$ws.Range("N113").Value = -6824
$ws.Range("H131").Value = 5209147.5
$ws.Range("J131").Value = 832.9892599999999
$ws.Range("L131").Value = 2498.96778
$ws.Range("N131").Value = -12578.96778
$ws.Range("H139").Value = 1996.6
$ws.Range("I139").Value = 1244.9
$ws.Range("J139").Value = 3500
$ws.Range("K139").Value = 3734.7
$ws.Range("L139").Value = 10500
$ws.Range("M139").Value = 1405.3
$ws.Range("N139").Value = -20780

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H132").Value = 2476.182
$ws.Range("I132").Value = 1900.4117
$ws.Range("J132").Value = 3087.9375
$ws.Range("K132").Value = 5701.2351
$ws.Range("L132").Value = 9263.8125
$ws.Range("M132").Value = -3171.2351
$ws.Range("N132").Value = -14323.8125

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2827.9285
$ws.Range("I136").Value = 1116.2174
$ws.Range("J136").Value = 4900
$ws.Range("K136").Value = 3348.6522
$ws.Range("L136").Value = 14700
$ws.Range("M136").Value = -798.6522
$ws.Range("N136").Value = -19800

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H132").Value = 12822369
$ws.Range("I132").Value = 1057.2727
$ws.Range("J132").Value = 83339580
$ws.Range("K132").Value = 3171.8181
$ws.Range("L132").Value = 250018740
$ws.Range("M132").Value = -641.8181
$ws.Range("N132").Value = -250023800
$ws.Range("H136").Value = 5000.769
$ws.Range("I136").Value = 3842.2
$ws.Range("J136").Value = 5724.875
$ws.Range("K136").Value = 11526.6
$ws.Range("L136").Value = 17174.625
$ws.Range("M136").Value = -8976.599999999999
$ws.Range("N136").Value = -22274.625
